# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in the Price (D) column whose new values look like plain decimal
# numbers must be pre-formatted as Text, otherwise Excel COM auto-converts
# the assigned string into a numeric value (losing the original text-cell type).
$dTextFixCells = @("D5", "D6", "D8", "D15", "D17", "D18", "D19", "D23", "D24", "D25", "D31", "D42", "D44", "D46", "D48", "D50")
foreach ($cellRef in $dTextFixCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "26.787.64"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.564.78"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "206.27"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "21.95"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "1.786.98"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.564.13"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "0.514"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.804.02"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "61.49"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "214.07"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "0.0₃0675"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "9.38"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "152.92"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "1.383.86"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "0.994"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").Value = "1.78"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "63.31"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "1.700.30"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "85.58"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "0.0₇0983"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "0.0950"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("E51").Value = "  -0.67%  "

# Restore default (Normal) style on the cells we temporarily reformatted,
# so no stray style index is introduced into the saved workbook.
foreach ($cellRef in $dTextFixCells) {
    $ws.Range($cellRef).Style = "Normal"
}
